$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (mushrooms)
$ws.Range("A2").Value = 287.2325069904327
$ws.Range("B2").Value = 376.6272072792053
$ws.Range("C2").Value = 54.32679343223572
$ws.Range("D2").Value = 114.1173338890076

# Row 3 (adult)
$ws.Range("A3").Value = 798.7408525943756
$ws.Range("B3").Value = 952.598132610321
$ws.Range("C3").Value = 71.91868281364441
$ws.Range("D3").Value = 116.7516672611237

# Row 4 (churn)
$ws.Range("A4").Value = 205.6780183315277
$ws.Range("B4").Value = 360.6242415904999
$ws.Range("C4").Value = 56.88042807579041
$ws.Range("D4").Value = 102.6660044193268

# Row 5 (credit card)
$ws.Range("A5").Value = 477.3747253417969
$ws.Range("B5").Value = 284.0919981002808
$ws.Range("C5").Value = 40.16576623916626
$ws.Range("D5").Value = 107.3109092712402
